$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the SQL queries stored in column B (TabQuery) and C2 (StatQuery).
#    The join conditions were changed from the generic ".id" columns to the
#    fully-qualified "<table>_id" columns, across every query cell.
# ---------------------------------------------------------------------------
$cells = @(
    $ws.Cells.Item(2,3),  # C2 - StatQuery
    $ws.Cells.Item(2,2),  # B2 - TabQuery (StudiesTab)
    $ws.Cells.Item(3,2),  # B3 - TabQuery (ParticipantsTab)
    $ws.Cells.Item(4,2),  # B4 - TabQuery (DiagnosisTab)
    $ws.Cells.Item(5,2),  # B5 - TabQuery (TreatmentTab)
    $ws.Cells.Item(6,2),  # B6 - TabQuery (TreatmentRespTab)
    $ws.Cells.Item(7,2)   # B7 - TabQuery (SurvivalTab)
)

foreach ($cell in $cells) {
    $text = $cell.Value2

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $cell.Value2 = $text
}

# The very first query (C2) also gained an extra trailing space after the
# final "WHERE" keyword (only this one query - the others are untouched).
$c2 = $ws.Cells.Item(2,3)
$newline = [char]10
$oldWhereLine = 'WHERE ' + $newline + '   std.dbgap_accession'
$newWhereLine = 'WHERE  ' + $newline + '   std.dbgap_accession'
$c2.Value2 = $c2.Value2.Replace($oldWhereLine, $newWhereLine)

# ---------------------------------------------------------------------------
# 2. Widen column C and drop its "best fit" flag (now a fixed custom width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 70.8333333333333

# ---------------------------------------------------------------------------
# 3. Move the active selection from C7 to C4 (scrolls the view back to the
#    top, clearing the old topLeftCell="A6" pin).
# ---------------------------------------------------------------------------
$null = $ws.Range("C4").Select()
